# Updates the cryptos list (price + 1h volume change) for Mon Apr 3 2023 run.
# Row 39/40 swap order (Aptos now ranks above InternetComputer) and all B/C/D/E
# values for that pair are replaced accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.975.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4974"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +17.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.098"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.011"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.813.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.231"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06585"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.973"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.015.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.218"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.024.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.367"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.044"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.632"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06842"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.914"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02305"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2135"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.928"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.148"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.291"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.666"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.940"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.172"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06751"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
